# Auto-generated Excel COM-interop script to apply market-data value updates
# across all 8 crafting-class worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 600.2857
$ws.Range("I2").Value = 149.25
$ws.Range("J2").Value = 780.7
$ws.Range("K2").Value = 149.25
$ws.Range("L2").Value = 780.7
$ws.Range("M2").Value = -36.25
$ws.Range("N2").Value = -1006.7
$ws.Range("H12").Value = 990
$ws.Range("J12").Value = 3500
$ws.Range("L12").Value = 3500
$ws.Range("N12").Value = -3840
$ws.Range("H15").Value = 1297.3662
$ws.Range("I15").Value = 1297.3662
$ws.Range("K15").Value = 3892.0986
$ws.Range("M15").Value = -3723.0986
$ws.Range("H17").Value = 552730.7
$ws.Range("J17").Value = 567938.5
$ws.Range("L17").Value = 1703815.5
$ws.Range("N17").Value = -1704151.5
$ws.Range("H19").Value = 1198.2222
$ws.Range("I19").Value = 1181.5
$ws.Range("J19").Value = 1231.6666
$ws.Range("K19").Value = 1181.5
$ws.Range("L19").Value = 1231.6666
$ws.Range("M19").Value = -1006.5
$ws.Range("N19").Value = -1581.6666
$ws.Range("H38").Value = 568.1667
$ws.Range("I38").Value = 560
$ws.Range("K38").Value = 1680
$ws.Range("M38").Value = -1308
$ws.Range("H40").Value = 84401.24000000001
$ws.Range("I40").Value = 376456.25
$ws.Range("K40").Value = 376456.25
$ws.Range("M40").Value = -376281.25
$ws.Range("H51").Value = 3497.0435
$ws.Range("J51").Value = 1943.75
$ws.Range("L51").Value = 1943.75
$ws.Range("N51").Value = -2911.75
$ws.Range("H53").Value = 830.0454999999999
$ws.Range("I53").Value = 565.4666999999999
$ws.Range("J53").Value = 1397
$ws.Range("K53").Value = 565.4666999999999
$ws.Range("L53").Value = 1397
$ws.Range("M53").Value = 71.53330000000005
$ws.Range("N53").Value = -2671
$ws.Range("H58").Value = 1719.2727
$ws.Range("I58").Value = 416
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 1248
$ws.Range("L58").Value = 12000
$ws.Range("M58").Value = -1098
$ws.Range("N58").Value = -12300
$ws.Range("H70").Value = 1891.2667
$ws.Range("J70").Value = 1990.6923
$ws.Range("L70").Value = 5972.0769
$ws.Range("N70").Value = -6512.0769
$ws.Range("H73").Value = 1891.2667
$ws.Range("J73").Value = 1990.6923
$ws.Range("L73").Value = 5972.0769
$ws.Range("N73").Value = -7844.0769
$ws.Range("H74").Value = 7791.5
$ws.Range("I74").Value = 10216.4
$ws.Range("K74").Value = 10216.4
$ws.Range("M74").Value = -9280.4
$ws.Range("H77").Value = 7791.5
$ws.Range("I77").Value = 10216.4
$ws.Range("K77").Value = 51082
$ws.Range("M77").Value = -46402
$ws.Range("H94").Value = 617.375
$ws.Range("I94").Value = 617.375
$ws.Range("K94").Value = 617.375
$ws.Range("M94").Value = -166.375
$ws.Range("H98").Value = 9199.777
$ws.Range("I98").Value = 9114.429
$ws.Range("J98").Value = 9498.5
$ws.Range("K98").Value = 9114.429
$ws.Range("L98").Value = 9498.5
$ws.Range("M98").Value = -7616.429
$ws.Range("N98").Value = -12494.5
$ws.Range("H104").Value = 1328.6
$ws.Range("I104").Value = 828
$ws.Range("J104").Value = 1453.75
$ws.Range("K104").Value = 2484
$ws.Range("L104").Value = 4361.25
$ws.Range("M104").Value = -737
$ws.Range("N104").Value = -7855.25
$ws.Range("H111").Value = 1587.5
$ws.Range("I111").Value = 1173
$ws.Range("J111").Value = 1836.2
$ws.Range("K111").Value = 3519
$ws.Range("L111").Value = 5508.6
$ws.Range("M111").Value = -452
$ws.Range("N111").Value = -11642.6
$ws.Range("H116").Value = 100032500
$ws.Range("I116").Value = 100032500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 100032500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -100029058
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 9199.777
$ws.Range("I122").Value = 9114.429
$ws.Range("J122").Value = 9498.5
$ws.Range("K122").Value = 27343.287
$ws.Range("L122").Value = 28495.5
$ws.Range("M122").Value = -24893.287
$ws.Range("N122").Value = -33395.5
$ws.Range("H132").Value = 2490.88
$ws.Range("I132").Value = 2511.3333
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 7533.999899999999
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -5003.999899999999
$ws.Range("N132").Value = -11060
$ws.Range("H135").Value = 1972.2142
$ws.Range("I135").Value = 1146.2
$ws.Range("K135").Value = 10315.8
$ws.Range("M135").Value = -7780.800000000001
$ws.Range("H137").Value = 1772.4546
$ws.Range("I137").Value = 1624.75
$ws.Range("J137").Value = 1856.8572
$ws.Range("K137").Value = 4874.25
$ws.Range("L137").Value = 5570.571599999999
$ws.Range("M137").Value = -2324.25
$ws.Range("N137").Value = -10670.5716
$ws.Range("H138").Value = 4171.96
$ws.Range("I138").Value = 2816.2
$ws.Range("J138").Value = 4411.212
$ws.Range("K138").Value = 8448.599999999999
$ws.Range("L138").Value = 13233.636
$ws.Range("M138").Value = -3308.599999999999
$ws.Range("N138").Value = -23513.636
$ws.Range("H141").Value = 4825.737
$ws.Range("I141").Value = 4162.3125
$ws.Range("K141").Value = 12486.9375
$ws.Range("M141").Value = -7306.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21145298
$ws.Range("I32").Value = 20318204
$ws.Range("K32").Value = 20318204
$ws.Range("M32").Value = -20317917
$ws.Range("H61").Value = 3333.1562
$ws.Range("I61").Value = 2988.1
$ws.Range("K61").Value = 2988.1
$ws.Range("M61").Value = -2776.1
$ws.Range("H63").Value = 6900
$ws.Range("I63").Value = 2600
$ws.Range("J63").Value = 8333.333000000001
$ws.Range("K63").Value = 2600
$ws.Range("L63").Value = 8333.333000000001
$ws.Range("M63").Value = -1914
$ws.Range("N63").Value = -9705.333000000001
$ws.Range("H64").Value = 30000
$ws.Range("I64").Value = 30000
$ws.Range("K64").Value = 30000
$ws.Range("M64").Value = -29752
$ws.Range("H66").Value = 6900
$ws.Range("I66").Value = 2600
$ws.Range("J66").Value = 8333.333000000001
$ws.Range("K66").Value = 13000
$ws.Range("L66").Value = 41666.665
$ws.Range("M66").Value = -9568
$ws.Range("N66").Value = -48530.665
$ws.Range("H67").Value = 30000
$ws.Range("I67").Value = 30000
$ws.Range("K67").Value = 30000
$ws.Range("M67").Value = -29142
$ws.Range("H74").Value = 1466.4445
$ws.Range("I74").Value = 1406.5807
$ws.Range("K74").Value = 1406.5807
$ws.Range("M74").Value = -532.5807
$ws.Range("H77").Value = 1466.4445
$ws.Range("I77").Value = 1406.5807
$ws.Range("K77").Value = 7032.9035
$ws.Range("M77").Value = -2664.9035
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 71965.664
$ws.Range("J93").Value = 71965.664
$ws.Range("L93").Value = 71965.664
$ws.Range("N93").Value = -76957.664
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 2407.4614
$ws.Range("I97").Value = 1542.5714
$ws.Range("J97").Value = 3416.5
$ws.Range("K97").Value = 1542.5714
$ws.Range("L97").Value = 3416.5
$ws.Range("M97").Value = -1046.5714
$ws.Range("N97").Value = -4408.5
$ws.Range("H102").Value = 2111
$ws.Range("I102").Value = 1527.75
$ws.Range("K102").Value = 1527.75
$ws.Range("M102").Value = 94.25
$ws.Range("H110").Value = 1678.5
$ws.Range("I110").Value = 943.7857
$ws.Range("K110").Value = 943.7857
$ws.Range("M110").Value = 1101.2143
$ws.Range("H122").Value = 2663.3635
$ws.Range("I122").Value = 2365.7778
$ws.Range("K122").Value = 7097.3334
$ws.Range("M122").Value = -4647.3334
$ws.Range("H132").Value = 671998.1
$ws.Range("I132").Value = 1003997.5
$ws.Range("J132").Value = 7999.4
$ws.Range("K132").Value = 3011992.5
$ws.Range("L132").Value = 23998.2
$ws.Range("M132").Value = -3009462.5
$ws.Range("N132").Value = -29058.2
$ws.Range("H133").Value = 53333.332
$ws.Range("J133").Value = 53333.332
$ws.Range("L133").Value = 53333.332
$ws.Range("N133").Value = -58393.332
$ws.Range("H136").Value = 3333.1562
$ws.Range("I136").Value = 2988.1
$ws.Range("K136").Value = 8964.299999999999
$ws.Range("M136").Value = -6414.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 18685.334
$ws.Range("J9").Value = 18685.334
$ws.Range("L9").Value = 18685.334
$ws.Range("N9").Value = -19021.334
$ws.Range("H81").Value = 70163.5
$ws.Range("J81").Value = 70163.5
$ws.Range("L81").Value = 70163.5
$ws.Range("N81").Value = -72285.5
$ws.Range("H82").Value = 28747.867
$ws.Range("I82").Value = 17862.924
$ws.Range("J82").Value = 99500
$ws.Range("K82").Value = 17862.924
$ws.Range("L82").Value = 99500
$ws.Range("M82").Value = -17479.924
$ws.Range("N82").Value = -100266
$ws.Range("H84").Value = 70163.5
$ws.Range("J84").Value = 70163.5
$ws.Range("L84").Value = 210490.5
$ws.Range("N84").Value = -221098.5
$ws.Range("H85").Value = 28747.867
$ws.Range("I85").Value = 17862.924
$ws.Range("J85").Value = 99500
$ws.Range("K85").Value = 17862.924
$ws.Range("L85").Value = 99500
$ws.Range("M85").Value = -16536.924
$ws.Range("N85").Value = -102152
$ws.Range("H94").Value = 459.16666
$ws.Range("I94").Value = 451.1
$ws.Range("K94").Value = 451.1
$ws.Range("M94").Value = -0.1000000000000227
$ws.Range("H134").Value = 2086954.9
$ws.Range("I134").Value = 2669649.5
$ws.Range("K134").Value = 8008948.5
$ws.Range("M134").Value = -8006413.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 81491.75
$ws.Range("J18").Value = 91989
$ws.Range("L18").Value = 91989
$ws.Range("N18").Value = -92449
$ws.Range("H31").Value = 2111.1829
$ws.Range("I31").Value = 1236.4791
$ws.Range("J31").Value = 3044.2
$ws.Range("K31").Value = 1236.4791
$ws.Range("L31").Value = 3044.2
$ws.Range("M31").Value = -941.4791
$ws.Range("N31").Value = -3634.2
$ws.Range("H34").Value = 2111.1829
$ws.Range("I34").Value = 1236.4791
$ws.Range("J34").Value = 3044.2
$ws.Range("K34").Value = 1236.4791
$ws.Range("L34").Value = 3044.2
$ws.Range("M34").Value = -1034.4791
$ws.Range("N34").Value = -3448.2
$ws.Range("H58").Value = 3385.9592
$ws.Range("I58").Value = 2701.25
$ws.Range("J58").Value = 4298.905
$ws.Range("K58").Value = 2701.25
$ws.Range("L58").Value = 4298.905
$ws.Range("M58").Value = -2498.25
$ws.Range("N58").Value = -4704.905
$ws.Range("H87").Value = 65370.8
$ws.Range("J87").Value = 65370.8
$ws.Range("L87").Value = 65370.8
$ws.Range("N87").Value = -67742.8
$ws.Range("H90").Value = 65370.8
$ws.Range("J90").Value = 65370.8
$ws.Range("L90").Value = 196112.4
$ws.Range("N90").Value = -207968.4
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H94").Value = 1131.2941
$ws.Range("I94").Value = 687.8
$ws.Range("J94").Value = 1316.0834
$ws.Range("K94").Value = 687.8
$ws.Range("L94").Value = 1316.0834
$ws.Range("M94").Value = -236.8
$ws.Range("N94").Value = -2218.0834
$ws.Range("H99").Value = 2852.1765
$ws.Range("I99").Value = 2689.5557
$ws.Range("J99").Value = 3035.125
$ws.Range("K99").Value = 2689.5557
$ws.Range("L99").Value = 3035.125
$ws.Range("M99").Value = -1191.5557
$ws.Range("N99").Value = -6031.125
$ws.Range("H103").Value = 60845.777
$ws.Range("J103").Value = 86660.336
$ws.Range("L103").Value = 86660.336
$ws.Range("N103").Value = -89004.336
$ws.Range("H106").Value = 70832.5
$ws.Range("I106").Value = 30000
$ws.Range("J106").Value = 84443.336
$ws.Range("K106").Value = 30000
$ws.Range("L106").Value = 84443.336
$ws.Range("M106").Value = -28738
$ws.Range("N106").Value = -86967.336
$ws.Range("H126").Value = 2852.1765
$ws.Range("I126").Value = 2689.5557
$ws.Range("J126").Value = 3035.125
$ws.Range("K126").Value = 8068.6671
$ws.Range("L126").Value = 9105.375
$ws.Range("M126").Value = -5598.6671
$ws.Range("N126").Value = -14045.375
$ws.Range("H132").Value = 2912.3914
$ws.Range("I132").Value = 2806.641
$ws.Range("J132").Value = 3501.5715
$ws.Range("K132").Value = 8419.923000000001
$ws.Range("L132").Value = 10504.7145
$ws.Range("M132").Value = -5889.923000000001
$ws.Range("N132").Value = -15564.7145
$ws.Range("H134").Value = 3267.2927
$ws.Range("I134").Value = 3229.4614
$ws.Range("K134").Value = 9688.3842
$ws.Range("M134").Value = -7153.3842
$ws.Range("H136").Value = 3385.9592
$ws.Range("I136").Value = 2701.25
$ws.Range("J136").Value = 4298.905
$ws.Range("K136").Value = 8103.75
$ws.Range("L136").Value = 12896.715
$ws.Range("M136").Value = -5553.75
$ws.Range("N136").Value = -17996.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 152196290
$ws.Range("I4").Value = 123621416
$ws.Range("J4").Value = 266495730
$ws.Range("K4").Value = 370864248
$ws.Range("L4").Value = 799487190
$ws.Range("M4").Value = -370864136
$ws.Range("N4").Value = -799487414
$ws.Range("H39").Value = 4229.6924
$ws.Range("J39").Value = 4229.6924
$ws.Range("L39").Value = 12689.0772
$ws.Range("N39").Value = -13277.0772
$ws.Range("H64").Value = 650
$ws.Range("I64").Value = 650
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 1950
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -1680
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 650
$ws.Range("I67").Value = 650
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 1950
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -1014
$ws.Range("N67").ClearContents()
$ws.Range("H68").Value = 2344.5557
$ws.Range("I68").Value = 1866.6666
$ws.Range("J68").Value = 2583.5
$ws.Range("K68").Value = 5599.9998
$ws.Range("L68").Value = 7750.5
$ws.Range("M68").Value = -4788.9998
$ws.Range("N68").Value = -9372.5
$ws.Range("H71").Value = 2344.5557
$ws.Range("I71").Value = 1866.6666
$ws.Range("J71").Value = 2583.5
$ws.Range("K71").Value = 16799.9994
$ws.Range("L71").Value = 23251.5
$ws.Range("M71").Value = -12743.9994
$ws.Range("N71").Value = -31363.5
$ws.Range("H107").Value = 1219.8
$ws.Range("I107").Value = 1338.8
$ws.Range("J107").Value = 1100.8
$ws.Range("K107").Value = 4016.4
$ws.Range("L107").Value = 3302.4
$ws.Range("M107").Value = -2096.4
$ws.Range("N107").Value = -7142.4
$ws.Range("H119").Value = 2753.1333
$ws.Range("I119").Value = 2191.4167
$ws.Range("K119").Value = 6574.250100000001
$ws.Range("M119").Value = -1736.250100000001
$ws.Range("H122").Value = 418.25
$ws.Range("J122").Value = 512.375
$ws.Range("L122").Value = 4611.375
$ws.Range("N122").Value = -9511.375
$ws.Range("H131").Value = 1815.4242
$ws.Range("I131").Value = 2001.9231
$ws.Range("J131").Value = 1694.2
$ws.Range("K131").Value = 6005.7693
$ws.Range("L131").Value = 5082.6
$ws.Range("M131").Value = -965.7692999999999
$ws.Range("N131").Value = -15162.6
$ws.Range("H140").Value = 1139.3572
$ws.Range("I140").Value = 995.1
$ws.Range("K140").Value = 2985.3
$ws.Range("M140").Value = 2194.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 268.09525
$ws.Range("I2").Value = 286.33334
$ws.Range("K2").Value = 286.33334
$ws.Range("M2").Value = -173.33334
$ws.Range("H80").Value = 3167.6667
$ws.Range("J80").Value = 3121.2
$ws.Range("L80").Value = 3121.2
$ws.Range("N80").Value = -5117.2
$ws.Range("H83").Value = 3167.6667
$ws.Range("J83").Value = 3121.2
$ws.Range("L83").Value = 15606
$ws.Range("N83").Value = -25590
$ws.Range("H102").Value = 2482.3572
$ws.Range("I102").Value = 2300.68
$ws.Range("K102").Value = 2300.68
$ws.Range("M102").Value = -678.6799999999998
$ws.Range("H113").Value = 25500.107
$ws.Range("I113").Value = 4489.9375
$ws.Range("J113").Value = 41507.855
$ws.Range("K113").Value = 4489.9375
$ws.Range("L113").Value = 41507.855
$ws.Range("M113").Value = -2319.9375
$ws.Range("N113").Value = -45847.855
$ws.Range("H132").Value = 3558.7144
$ws.Range("I132").Value = 3735.889
$ws.Range("J132").Value = 3239.8
$ws.Range("K132").Value = 11207.667
$ws.Range("L132").Value = 9719.400000000001
$ws.Range("M132").Value = -8677.667000000001
$ws.Range("N132").Value = -14779.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value = 14496291
$ws.Range("I40").Value = 15876558
$ws.Range("K40").Value = 15876558
$ws.Range("M40").Value = -15876422
$ws.Range("H46").Value = 3096.2632
$ws.Range("I46").Value = 1980.2
$ws.Range("J46").Value = 3265.3635
$ws.Range("K46").Value = 1980.2
$ws.Range("L46").Value = 3265.3635
$ws.Range("M46").Value = -1792.2
$ws.Range("N46").Value = -3641.3635
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 593585.5
$ws.Range("I132").Value = 593585.5
$ws.Range("K132").Value = 1780756.5
$ws.Range("M132").Value = -1778226.5
$ws.Range("H136").Value = 41219
$ws.Range("I136").Value = 56312.145
$ws.Range("K136").Value = 168936.435
$ws.Range("M136").Value = -166386.435
$ws.Range("H139").Value = 100885
$ws.Range("I139").Value = 70000
$ws.Range("J139").Value = 121475
$ws.Range("K139").Value = 70000
$ws.Range("L139").Value = 121475
$ws.Range("M139").Value = -64860
$ws.Range("N139").Value = -131755

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 19638.818
$ws.Range("J15").Value = 19638.818
$ws.Range("L15").Value = 19638.818
$ws.Range("N15").Value = -20214.818
$ws.Range("H74").Value = 8399.799999999999
$ws.Range("J74").Value = 8399.799999999999
$ws.Range("L74").Value = 8399.799999999999
$ws.Range("N74").Value = -10271.8
$ws.Range("H77").Value = 8399.799999999999
$ws.Range("J77").Value = 8399.799999999999
$ws.Range("L77").Value = 25199.4
$ws.Range("N77").Value = -34559.39999999999
$ws.Range("H122").Value = 55562400
$ws.Range("J122").Value = 7224.5
$ws.Range("L122").Value = 21673.5
$ws.Range("N122").Value = -26573.5
$ws.Range("H132").Value = 26592.715
$ws.Range("I132").Value = 29122.658
$ws.Range("J132").Value = 2558.25
$ws.Range("K132").Value = 87367.974
$ws.Range("L132").Value = 7674.75
$ws.Range("M132").Value = -84837.974
$ws.Range("N132").Value = -12734.75
$ws.Range("H135").Value = 149999.25
$ws.Range("J135").Value = 149999.25
$ws.Range("L135").Value = 149999.25
$ws.Range("N135").Value = -160139.25
$ws.Range("H136").Value = 32175.686
$ws.Range("I136").Value = 3238
$ws.Range("J136").Value = 81147.16
$ws.Range("K136").Value = 9714
$ws.Range("L136").Value = 243441.48
$ws.Range("M136").Value = -7164
$ws.Range("N136").Value = -248541.48
